# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the display order of "Santa Lucia" (row 207) and "Nueva Caledonia" (row 208)
$ws.Range("A207").Value = "Nueva Caledonia"
$ws.Range("A208").Value = "Santa Lucia"

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 7 de Octubre de 2020 a las 08:31"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 7723721
$ws.Range("C4").Value = 975
$ws.Range("E4").Value = 2572340
$ws.Range("G4").Value = 14
$ws.Range("H4").Value = 215836

# Row 27 - Israel
$ws.Range("B27").Value = 278585
$ws.Range("C27").Value = 1559
$ws.Range("D27").Value = 215181
$ws.Range("E27").Value = 61601
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 1803

# Row 28 - Ucrania
$ws.Range("B28").Value = 239337
$ws.Range("C28").Value = 4753
$ws.Range("D28").Value = 105970
$ws.Range("E28").Value = 128770
$ws.Range("G28").Value = 77
$ws.Range("H28").Value = 4597

# Row 59 - Uzbekistan
$ws.Range("B59").Value = 59579
$ws.Range("C59").Value = 236
$ws.Range("D59").Value = 56165
$ws.Range("E59").Value = 2923
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = 491

# Row 79 - El Salvador
$ws.Range("E79").Value = 4236
$ws.Range("G79").Value = 4
$ws.Range("H79").Value = 873
